$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @([double]"0.06366858631372452", [double]"0.9831519722938538", [double]"0.008906668983399868", [double]"0.9985330104827881"),
    @([double]"0.01061982475221157", [double]"0.9983243942260742", [double]"0.005267607979476452", [double]"0.9987286329269409"),
    @([double]"0.005555095616728067", [double]"0.9988215565681458", [double]"0.001552314730361104", [double]"0.9993153810501099"),
    @([double]"0.003641140414401889", [double]"0.9992450475692749", [double]"0.0004941453225910664", [double]"0.9998043775558472"),
    @([double]"0.002136968309059739", [double]"0.9995765089988708", [double]"0.0002020417159656063", [double]"1"),
    @([double]"0.001222585560753942", [double]"0.999760627746582", [double]"0.0001328032521996647", [double]"1"),
    @([double]"0.0009182837675325572", [double]"0.9997237920761108", [double]"9.245761611964554E-05", [double]"1"),
    @([double]"0.001980241155251861", [double]"0.9996317625045776", [double]"0.0001364320487482473", [double]"1"),
    @([double]"0.0009803419234231114", [double]"0.9997790455818176", [double]"0.0001039370545186102", [double]"1"),
    @([double]"0.0006622901419177651", [double]"0.9997974634170532", [double]"8.127058390527964E-05", [double]"1"),
    @([double]"0.001301989075727761", [double]"0.9997053742408752", [double]"4.765586345456541E-05", [double]"1"),
    @([double]"0.0005972832441329956", [double]"0.9999263286590576", [double]"7.329156505875289E-05", [double]"1"),
    @([double]"0.0005894917412661016", [double]"0.9998158812522888", [double]"1.234308001585305E-05", [double]"1"),
    @([double]"0.001299133407883346", [double]"0.9997790455818176", [double]"5.894942205486586E-06", [double]"1"),
    @([double]"0.0004160820972174406", [double]"0.9998711347579956", [double]"3.818998720817035E-06", [double]"1"),
    @([double]"0.0003255673218518496", [double]"0.9998894929885864", [double]"9.798271776162437E-07", [double]"1"),
    @([double]"0.0006824019947089255", [double]"0.999760627746582", [double]"2.34931540035177E-05", [double]"1"),
    @([double]"7.379411545116454E-05", [double]"0.9999815821647644", [double]"0.000165854042279534", [double]"1"),
    @([double]"0.001373500796034932", [double]"0.9997790455818176", [double]"1.569732717143779E-06", [double]"1"),
    @([double]"0.0006420322461053729", [double]"0.9997790455818176", [double]"1.524807157693431E-05", [double]"1"),
    @([double]"0.0006171134882606566", [double]"0.9998711347579956", [double]"2.487443452992011E-05", [double]"1"),
    @([double]"0.0005082093412056565", [double]"0.99985271692276", [double]"0.0007529320428147912", [double]"0.9998043775558472"),
    @([double]"0.0002891742915380746", [double]"0.9999263286590576", [double]"5.739016614825232E-06", [double]"1"),
    @([double]"0.0001173531345557421", [double]"0.9999815821647644", [double]"1.989830889215227E-05", [double]"1"),
    @([double]"9.582204802427441E-05", [double]"0.9999447464942932", [double]"5.464713922265219E-06", [double]"1"),
    @([double]"0.0008906829752959311", [double]"0.999907910823822", [double]"1.481008240489246E-07", [double]"1"),
    @([double]"0.0002216367283836007", [double]"0.999907910823822", [double]"5.635368438561272E-08", [double]"1"),
    @([double]"4.225752491038293E-05", [double]"1", [double]"3.041580853846426E-08", [double]"1"),
    @([double]"0.0003786913875956088", [double]"0.9999263286590576", [double]"2.248880264232866E-06", [double]"1"),
    @([double]"0.0002162736491300166", [double]"0.9999263286590576", [double]"3.182578112159717E-08", [double]"1"),
    @([double]"0.0003994805156253278", [double]"0.9999631643295288", [double]"1.375420595195465E-07", [double]"1"),
    @([double]"0.0002254263672512025", [double]"0.9999263286590576", [double]"5.171299903850013E-07", [double]"1"),
    @([double]"0.0003270395682193339", [double]"0.9999447464942932", [double]"3.800854756264016E-06", [double]"1"),
    @([double]"0.0007396579021587968", [double]"0.9998711347579956", [double]"1.347218585578958E-07", [double]"1"),
    @([double]"0.0001641543203731999", [double]"0.9999447464942932", [double]"2.171018991248275E-07", [double]"1"),
    @([double]"7.046959945000708E-05", [double]"0.9999815821647644", [double]"5.181912001717137E-07", [double]"1"),
    @([double]"4.265984534868039E-05", [double]"0.9999815821647644", [double]"8.827152669255156E-06", [double]"1"),
    @([double]"1.936934131663293E-05", [double]"0.9999815821647644", [double]"3.473567858236493E-06", [double]"1"),
    @([double]"0.0005887700244784355", [double]"0.9998894929885864", [double]"0.0002955517265945673", [double]"0.9999021887779236"),
    @([double]"0.0003093911800533533", [double]"0.9999447464942932", [double]"4.290322763722543E-09", [double]"1"),
    @([double]"0.000224388757487759", [double]"0.9999263286590576", [double]"3.685969929279054E-08", [double]"1"),
    @([double]"0.0003018028510268778", [double]"0.9999631643295288", [double]"6.824495102364381E-08", [double]"1"),
    @([double]"9.885265171760693E-06", [double]"1", [double]"3.633481782117087E-08", [double]"1"),
    @([double]"0.0001811587135307491", [double]"0.9999263286590576", [double]"2.005269283955613E-09", [double]"1"),
    @([double]"0.0005346070975065231", [double]"0.999907910823822", [double]"4.253402039466891E-06", [double]"1"),
    @([double]"0.0002046288136625662", [double]"0.9999447464942932", [double]"1.447288298095373E-07", [double]"1"),
    @([double]"0.0003818174882326275", [double]"0.9998894929885864", [double]"3.536594590514142E-07", [double]"1"),
    @([double]"0.000231315178098157", [double]"0.9998894929885864", [double]"1.325510545768793E-08", [double]"1"),
    @([double]"2.351991861360148E-05", [double]"0.9999815821647644", [double]"6.085639547848132E-09", [double]"1"),
    @([double]"1.007928858598461E-05", [double]"1", [double]"1.375708436057721E-09", [double]"1")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}

Write-Host "Updated A2:D51 with new training history values (20 LSTMs, 50 dense units run)"
